$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.838.21"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "1.966.67"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "324.46"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.4790"
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  -4.78%  "
$ws.Range("D9").Value = "53.97"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "0.08588"
$ws.Range("E10").Value = "  -6.28%  "
$ws.Range("D11").Value = "1.063"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("D12").Value = "22.51"
$ws.Range("E12").Value = "  -4.41%  "
$ws.Range("D13").Value = "1.960.78"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "7.708"
$ws.Range("E14").Value = "  -5.32%  "
$ws.Range("D15").Value = "6.265"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("D16").Value = "1.015"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "90.60"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "0.00001069"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").Value = "0.06624"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "18.72"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "5.806"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("D23").Value = "28.864.25"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").Value = "11.58"
$ws.Range("E24").Value = "  -4.03%  "
$ws.Range("D25").Value = "2.294"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "2.220.33"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "154.09"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").Value = "6.014"
$ws.Range("E29").Value = "  -6.36%  "
$ws.Range("D30").Value = "2.169"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("D31").Value = "124.78"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").Value = "1.009"
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").Value = "0.09654"
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("D34").Value = "1.467"
$ws.Range("E34").Value = "  -7.40%  "
$ws.Range("D35").Value = "5.712"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "3.700"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").Value = "0.02355"
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").Value = "1.278"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").Value = "0.06264"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "8.816"
$ws.Range("E40").Value = "  -7.85%  "
$ws.Range("D41").Value = "0.6259"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").Value = "11.22"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "0.1920"
$ws.Range("E44").Value = "  -7.47%  "
$ws.Range("D45").Value = "1.334"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").Value = "0.5973"
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("D47").Value = "13.09"
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").Value = "2.089"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("D49").Value = "3.443"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("D50").Value = "0.00000000336"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").Value = "2.113"
$ws.Range("E51").Value = "  +5.37%  "
